$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.566.62"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "2.337.72"
$ws.Range("E3").Value = "  -2.31%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.93%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -8.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -3.40%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.19%  "
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.963"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.21%  "
$ws.Range("D15").Value = "2.692.42"
$ws.Range("E15").Value = "  -2.47%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.74%  "
$ws.Range("D17").Value = "2.334.21"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").Value = "44.571.80"
$ws.Range("E18").Value = "  -2.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +12.36%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.77%  "
$ws.Range("E21").Value = "  -2.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "255.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.03%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.84"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.11%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0939"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "36.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "165.04"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.34%  "
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.19"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.46%  "
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.66"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.91"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0349"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.35%  "
$ws.Range("D42").Value = "1.868.83"
$ws.Range("E42").Value = "  +13.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "94.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.68%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.76"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "68.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.224"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.17%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "82.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.19"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.73%  "
